# Revised Data for consistency
# The column header in C1 is shortened from "Preventative Health" to "Preventative"
# to be consistent with the other single-word motivation headers (Wellness, At Risk,
# Sick Role, Self Care).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "Preventative"
